$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.765.01"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3
$ws.Range("D3").Value = "3.391.49"
$ws.Range("E3").Value = "  -1.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").Value = "'407.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6
$ws.Range("D6").Value = "'127.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "

# Row 7
$ws.Range("D7").Value = "'0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.90%  "

# Row 8
$ws.Range("E8").Value = "  +0.27%  "

# Row 9
$ws.Range("D9").Value = "'0.704"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.98%  "

# Row 10
$ws.Range("D10").Value = "'0.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.18%  "

# Row 11
$ws.Range("D11").Value = "'41.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "

# Row 12
$ws.Range("D12").Value = "3.948.00"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "'0.140"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

# Row 14
$ws.Range("D14").Value = "'8.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'20.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.02%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000201"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.28%  "

# Row 17
$ws.Range("D17").Value = "3.379.74"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").Value = "'12.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "

# Row 19
$ws.Range("D19").Value = "'1.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.54%  "

# Row 20
$ws.Range("D20").Value = "61.786.06"
$ws.Range("E20").Value = "  -0.27%  "

# Row 21
$ws.Range("D21").Value = "'468.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.77%  "

# Row 22
$ws.Range("D22").Value = "'88.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "

# Row 23
$ws.Range("D23").Value = "'3.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

# Row 24
$ws.Range("D24").Value = "'12.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.24%  "

# Row 25
$ws.Range("D25").Value = "'3.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.94%  "

# Row 26
$ws.Range("D26").Value = "'32.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "

# Row 27
$ws.Range("D27").Value = "'8.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

# Row 28
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("D29").Value = "'7.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.20%  "

# Row 30
$ws.Range("D30").Value = "'2.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.47%  "

# Row 31
$ws.Range("D31").Value = "'11.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "

# Row 32
$ws.Range("D32").Value = "'0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.56%  "

# Row 33
$ws.Range("D33").Value = "'0.162"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.86%  "

# Row 34
$ws.Range("D34").Value = "'40.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.88%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("D36").Value = "'55.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.56%  "

# Row 37
$ws.Range("D37").Value = "'0.0478"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.46%  "

# Row 38
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("D39").Value = "'148.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.15%  "

# Row 40
$ws.Range("D40").Value = "'3.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.88%  "

# Row 41
$ws.Range("D41").Value = "'0.132"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.25%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.310"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

# Row 44
$ws.Range("D44").Value = "'2.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.46%  "

# Row 45
$ws.Range("D45").Value = "'2.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.72%  "

# Row 46
$ws.Range("D46").Value = "'4.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.61%  "

# Row 48
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'15.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.44%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'21.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.141"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.65%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.748.17"
$ws.Range("E51").Value = "  -1.02%  "
